$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 197
$ws.Range("I33").Value = 197
$ws.Range("K33").Value = 197
$ws.Range("M33").Value = 32

$ws.Range("H75").Value = 30314
$ws.Range("J75").Value = 30314
$ws.Range("L75").Value = 30314
$ws.Range("N75").Value = -32186

$ws.Range("H78").Value = 30314
$ws.Range("J78").Value = 30314
$ws.Range("L78").Value = 90942
$ws.Range("N78").Value = -100302

$ws.Range("H96").Value = 1208.25
$ws.Range("I96").Value = 40
$ws.Range("J96").Value = 1597.6666
$ws.Range("K96").Value = 120
$ws.Range("L96").Value = 4792.9998
$ws.Range("M96").Value = 1253
$ws.Range("N96").Value = -7538.9998

$ws.Range("H107").Value = 778.7692
$ws.Range("I107").Value = 822.4
$ws.Range("J107").Value = 633.3333
$ws.Range("K107").Value = 822.4
$ws.Range("L107").Value = 633.3333
$ws.Range("M107").Value = 1097.6
$ws.Range("N107").Value = -4473.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1010.5
$ws.Range("I2").Value = 1010.5
$ws.Range("K2").Value = 1010.5
$ws.Range("M2").Value = -897.5

$ws.Range("H102").Value = 6065.273
$ws.Range("I102").Value = 3999.2
$ws.Range("J102").Value = 7787
$ws.Range("K102").Value = 3999.2
$ws.Range("L102").Value = 7787
$ws.Range("M102").Value = -2377.2
$ws.Range("N102").Value = -11031

$ws.Range("H110").Value = 1693
$ws.Range("I110").Value = 1275.2858
$ws.Range("J110").Value = 2667.6667
$ws.Range("K110").Value = 1275.2858
$ws.Range("L110").Value = 2667.6667
$ws.Range("M110").Value = 769.7141999999999
$ws.Range("N110").Value = -6757.6667

$ws.Range("H116").Value = 1010.5
$ws.Range("I116").Value = 1010.5
$ws.Range("K116").Value = 1010.5
$ws.Range("M116").Value = 1283.5

$ws.Range("H132").Value = 3707.1765
$ws.Range("I132").Value = 4119.0713
$ws.Range("J132").Value = 1785
$ws.Range("K132").Value = 12357.2139
$ws.Range("L132").Value = 5355
$ws.Range("M132").Value = -9827.213899999999
$ws.Range("N132").Value = -10415

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1010.5
$ws.Range("I3").Value = 1010.5
$ws.Range("K3").Value = 1010.5
$ws.Range("M3").Value = -896.5

$ws.Range("H20").Value = 971
$ws.Range("I20").Value = 959.6
$ws.Range("J20").Value = 999.5
$ws.Range("K20").Value = 959.6
$ws.Range("L20").Value = 999.5
$ws.Range("M20").Value = -712.6
$ws.Range("N20").Value = -1493.5

$ws.Range("H86").Value = 6642.7144
$ws.Range("I86").Value = 1749.5
$ws.Range("J86").Value = 8600
$ws.Range("K86").Value = 1749.5
$ws.Range("L86").Value = 8600
$ws.Range("M86").Value = -626.5
$ws.Range("N86").Value = -10846

$ws.Range("H89").Value = 6642.7144
$ws.Range("I89").Value = 1749.5
$ws.Range("J89").Value = 8600
$ws.Range("K89").Value = 8747.5
$ws.Range("L89").Value = 43000
$ws.Range("M89").Value = -3131.5
$ws.Range("N89").Value = -54232

$ws.Range("H134").Value = 1559
$ws.Range("I134").Value = 1559
$ws.Range("K134").Value = 4677
$ws.Range("M134").Value = -2142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 110578.8
$ws.Range("I16").Value = 122782
$ws.Range("K16").Value = 122782
$ws.Range("M16").Value = -122495

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H105").Value = 420.66666
$ws.Range("I105").Value = 264.8
$ws.Range("K105").Value = 264.8
$ws.Range("M105").Value = 1482.2

$ws.Range("H113").Value = 110578.8
$ws.Range("I113").Value = 122782
$ws.Range("K113").Value = 122782
$ws.Range("M113").Value = -120612

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

$ws.Range("H120").Value = 1400
$ws.Range("I120").Value = 1400
$ws.Range("K120").Value = 4200
$ws.Range("M120").Value = 638

$ws.Range("H134").Value = 1524.2858
$ws.Range("I134").Value = 1524.2858
$ws.Range("K134").Value = 4572.857400000001
$ws.Range("M134").Value = 497.1425999999992

$ws.Range("H139").Value = 1800
$ws.Range("I139").Value = 1025
$ws.Range("K139").Value = 3075
$ws.Range("M139").Value = 2065

$ws.Range("H140").Value = 1928.25
$ws.Range("I140").Value = 1782.6428
$ws.Range("K140").Value = 5347.928400000001
$ws.Range("M140").Value = -167.9284000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3987.2222
$ws.Range("I126").Value = 3147.4285
$ws.Range("K126").Value = 9442.2855
$ws.Range("M126").Value = -6972.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H46").Value = 3507.0938
$ws.Range("I46").Value = 3689.889
$ws.Range("J46").Value = 3435.5652
$ws.Range("K46").Value = 3689.889
$ws.Range("L46").Value = 3435.5652
$ws.Range("M46").Value = -3501.889
$ws.Range("N46").Value = -3811.5652

$ws.Range("H75").Value = 16404.666
$ws.Range("I75").Value = 2107
$ws.Range("J75").Value = 45000
$ws.Range("K75").Value = 2107
$ws.Range("L75").Value = 45000
$ws.Range("M75").Value = -1171
$ws.Range("N75").Value = -46872

$ws.Range("H78").Value = 16404.666
$ws.Range("I78").Value = 2107
$ws.Range("J78").Value = 45000
$ws.Range("K78").Value = 6321
$ws.Range("L78").Value = 135000
$ws.Range("M78").Value = -1641
$ws.Range("N78").Value = -144360

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H136").Value = 2218
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 832.6667
$ws.Range("I81").Value = 839.4
$ws.Range("J81").Value = 799
$ws.Range("K81").Value = 1678.8
$ws.Range("L81").Value = 1598
$ws.Range("M81").Value = -617.8
$ws.Range("N81").Value = -3720

$ws.Range("H84").Value = 832.6667
$ws.Range("I84").Value = 839.4
$ws.Range("J84").Value = 799
$ws.Range("K84").Value = 8394
$ws.Range("L84").Value = 7990
$ws.Range("M84").Value = -3090
$ws.Range("N84").Value = -18598

$ws.Range("H100").Value = 1742.5714
$ws.Range("I100").Value = 1039.8
$ws.Range("K100").Value = 2079.6
$ws.Range("M100").Value = -1538.6

$ws.Range("H126").Value = 4638.8
$ws.Range("I126").Value = 3251.9092
$ws.Range("J126").Value = 6333.8887
$ws.Range("K126").Value = 9755.7276
$ws.Range("L126").Value = 19001.6661
$ws.Range("M126").Value = -7285.7276
$ws.Range("N126").Value = -23941.6661

$ws.Range("H132").Value = 2640
$ws.Range("I132").Value = 2640
$ws.Range("K132").Value = 7920
$ws.Range("M132").Value = -5390
